# Applies the "Analog / Digital" screen text additions (plus the first,
# still-unused placeholders for the upcoming PWM & Accelerometer screen) to
# the TouchGFX texts workbook, per the commit:
#   "Analog Digital Screen Working here. Some additions for the PWM
#    Accelerometer screen are included but just started none of it does
#    anything"

function Set-Text($ws, $addr, $val) {
    # Writes a normal (string or number) value into a cell.
    $ws.Range($addr).Value = $val
}

function Set-TextForceString($ws, $addr, $val) {
    # Some of the new values look like numbers (e.g. "0", "0.000", "00.00")
    # but must be stored as text, matching the source data.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

function Set-EmptyPlaceholder($ws, $addr) {
    # Touch the cell so an (empty) <c> entry is materialized for it, just
    # like the sibling cells already on the row.
    $ws.Range($addr).NumberFormat = "General"
}

$wb = $excel.ActiveWorkbook
$wsTypography  = $wb.Worksheets.Item("Typography")
$wsTranslation = $wb.Worksheets.Item("Translation")

# ---------------------------------------------------------------------------
# Typography sheet
# Columns: B=Typography Name, C=Font, D=Size, E=Bpp, F=Fallback Character,
#          G=Wildcard Characters, H=Widget Wildcard Characters,
#          I=Wildcard Ranges, J=Ellipsis Character
# ---------------------------------------------------------------------------

# Row 6 (Small) gains a "Widget Wildcard Characters" value
Set-Text $wsTypography "H6" "-., 0123456789"

# Row 7 (Iceland_45) gains "Wildcard Characters" and "Wildcard Ranges"
Set-Text $wsTypography "G7" "."
Set-Text $wsTypography "I7" "0-9"

# New row 8: Iceland_200
Set-Text $wsTypography "B8" "Iceland_200"
Set-Text $wsTypography "C8" "Iceland-Regular.ttf"
Set-Text $wsTypography "D8" 200
Set-Text $wsTypography "E8" 4
Set-Text $wsTypography "F8" "?"
Set-EmptyPlaceholder $wsTypography "G8"
Set-EmptyPlaceholder $wsTypography "H8"
Set-Text $wsTypography "I8" "0-9,A-F"
Set-EmptyPlaceholder $wsTypography "J8"
$wsTypography.Range("B8:J8").Style = "Normal"

# New row 9: Digital_Dream_25
Set-Text $wsTypography "B9" "Digital_Dream_25"
Set-Text $wsTypography "C9" "digital-dream.fat-skew.ttf"
Set-Text $wsTypography "D9" 25
Set-Text $wsTypography "E9" 4
Set-Text $wsTypography "F9" "-"
Set-Text $wsTypography "G9" "."
Set-EmptyPlaceholder $wsTypography "H9"
Set-Text $wsTypography "I9" "0-9"
Set-EmptyPlaceholder $wsTypography "J9"
$wsTypography.Range("B9:J9").Style = "Normal"

# New row 10: Digital_Dream_100
Set-Text $wsTypography "B10" "Digital_Dream_100"
Set-Text $wsTypography "C10" "digital-dream.fat-skew.ttf"
Set-Text $wsTypography "D10" 100
Set-Text $wsTypography "E10" 4
Set-Text $wsTypography "F10" "-"
Set-EmptyPlaceholder $wsTypography "G10"
Set-EmptyPlaceholder $wsTypography "H10"
Set-Text $wsTypography "I10" "0-9,A-F"
Set-EmptyPlaceholder $wsTypography "J10"
$wsTypography.Range("B10:J10").Style = "Normal"

# ---------------------------------------------------------------------------
# Translation sheet
# Columns: B=TEXT ID, C=TYPOGRAPHY NAME, D=ALIGNMENT, E=DIRECTION, F=GB
# ---------------------------------------------------------------------------

# Existing row 6 (SingleUseId7): text shortened from "Analog & Digital" to
# "Analog"
Set-Text $wsTranslation "F6" "Analog"

# New row 8
Set-Text $wsTranslation "B8" "SingleUseId9"
Set-Text $wsTranslation "C8" "Iceland_200"
Set-Text $wsTranslation "D8" "Left"
Set-Text $wsTranslation "E8" "LTR"
Set-Text $wsTranslation "F8" "<Digital_Read>"

# New row 9
Set-Text $wsTranslation "B9" "SingleUseId10"
Set-Text $wsTranslation "C9" "Digital_Dream_25"
Set-Text $wsTranslation "D9" "Center"
Set-Text $wsTranslation "E9" "LTR"
Set-Text $wsTranslation "F9" "<Analog_Read>V"

# New row 10
Set-Text $wsTranslation "B10" "SingleUseId11"
Set-Text $wsTranslation "C10" "Iceland_45"
Set-Text $wsTranslation "D10" "Left"
Set-Text $wsTranslation "E10" "LTR"
Set-Text $wsTranslation "F10" "Digital"

# New row 11
Set-Text $wsTranslation "B11" "SingleUseId12"
Set-Text $wsTranslation "C11" "Digital_Dream_25"
Set-Text $wsTranslation "D11" "Left"
Set-Text $wsTranslation "E11" "LTR"
Set-TextForceString $wsTranslation "F11" "0.000"

# New row 12
Set-Text $wsTranslation "B12" "SingleUseId13"
Set-Text $wsTranslation "C12" "Iceland_200"
Set-Text $wsTranslation "D12" "Left"
Set-Text $wsTranslation "E12" "LTR"
Set-TextForceString $wsTranslation "F12" "0"

# New row 13
Set-Text $wsTranslation "B13" "SingleUseId16"
Set-Text $wsTranslation "C13" "Small"
Set-Text $wsTranslation "D13" "Left"
Set-Text $wsTranslation "E13" "LTR"
Set-Text $wsTranslation "F13" "<>"

# New row 14
Set-Text $wsTranslation "B14" "SingleUseId17"
Set-Text $wsTranslation "C14" "Small"
Set-Text $wsTranslation "D14" "Right"
Set-Text $wsTranslation "E14" "LTR"
Set-Text $wsTranslation "F14" "<>"

# New row 15
Set-Text $wsTranslation "B15" "SingleUseId18"
Set-Text $wsTranslation "C15" "Digital_Dream_25"
Set-Text $wsTranslation "D15" "Left"
Set-Text $wsTranslation "E15" "LTR"
Set-Text $wsTranslation "F15" "<value>"

# New row 16
Set-Text $wsTranslation "B16" "SingleUseId19"
Set-Text $wsTranslation "C16" "Digital_Dream_25"
Set-Text $wsTranslation "D16" "Left"
Set-Text $wsTranslation "E16" "LTR"
Set-TextForceString $wsTranslation "F16" "00.00"

# New row 17
Set-Text $wsTranslation "B17" "SingleUseId20"
Set-Text $wsTranslation "C17" "Digital_Dream_25"
Set-Text $wsTranslation "D17" "Left"
Set-Text $wsTranslation "E17" "LTR"
Set-Text $wsTranslation "F17" "<value>"

# New row 18
Set-Text $wsTranslation "B18" "SingleUseId21"
Set-Text $wsTranslation "C18" "Digital_Dream_25"
Set-Text $wsTranslation "D18" "Left"
Set-Text $wsTranslation "E18" "LTR"
Set-TextForceString $wsTranslation "F18" "00.00"

# New row 19
Set-Text $wsTranslation "B19" "SingleUseId22"
Set-Text $wsTranslation "C19" "Digital_Dream_25"
Set-Text $wsTranslation "D19" "Left"
Set-Text $wsTranslation "E19" "LTR"
Set-Text $wsTranslation "F19" "<value>"

# New row 20
Set-Text $wsTranslation "B20" "SingleUseId23"
Set-Text $wsTranslation "C20" "Digital_Dream_25"
Set-Text $wsTranslation "D20" "Left"
Set-Text $wsTranslation "E20" "LTR"
Set-TextForceString $wsTranslation "F20" "00.00"
